$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.730466365814209
$ws.Range("B1").Value = 1.818760514259338
$ws.Range("C1").Value = 1.98962140083313
$ws.Range("D1").Value = 2.843653917312622
$ws.Range("E1").Value = 5.152978897094727
